$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-7, columns A-H
$data = @(
    @(1, 0.5, 100000, 0.25, 0.25, 0.25, 0, 3),
    @(2, 0.5, 100000, 0.25, 0.25, 0.25, 0, 3),
    @(3, 1,   200000, 0.5,  0.5,  0.4,  0, 4),
    @(4, 1,   200000, 0.5,  0.5,  0.4,  0, 4),
    @(5, 2,   320000, 1,    1,    0.8,  0, 6),
    @(6, 2,   320000, 1,    1,    0.8,  0, 6)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Length; $c++) {
        $col = $c + 1
        $ws.Cells.Item($row, $col).Value = $rowData[$c]
    }
}

$ws.Range("C2").Select() | Out-Null
